$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 23333
$ws.Range("I21").Value = 17499.75
$ws.Range("J21").Value = 34999.5
$ws.Range("K21").Value = 17499.75
$ws.Range("L21").Value = 34999.5
$ws.Range("M21").Value = -17031.75
$ws.Range("N21").Value = -35935.5
$ws.Range("H23").Value = 23333
$ws.Range("I23").Value = 17499.75
$ws.Range("J23").Value = 34999.5
$ws.Range("K23").Value = 17499.75
$ws.Range("L23").Value = 34999.5
$ws.Range("M23").Value = -17265.75
$ws.Range("N23").Value = -35467.5
$ws.Range("H115").Value = 8181965.5
$ws.Range("J115").Value = 1302
$ws.Range("L115").Value = 3906
$ws.Range("N115").Value = -7040

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 2000333.4
$ws.Range("I6").Value = 2000333.4
$ws.Range("K6").Value = 2000333.4
$ws.Range("M6").Value = -2000160.4
$ws.Range("H15").Value = 880
$ws.Range("I15").Value = 750
$ws.Range("K15").Value = 750
$ws.Range("M15").Value = -400
$ws.Range("H16").Value = 3875
$ws.Range("I16").Value = 250
$ws.Range("K16").Value = 250
$ws.Range("M16").Value = 37
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()
$ws.Range("H97").Value = 3101420.5
$ws.Range("I97").Value = 6194459.5
$ws.Range("J97").Value = 8381.5
$ws.Range("K97").Value = 6194459.5
$ws.Range("L97").Value = 8381.5
$ws.Range("M97").Value = -6193963.5
$ws.Range("N97").Value = -9373.5
$ws.Range("H132").Value = 6903591
$ws.Range("J132").Value = 9219.214
$ws.Range("L132").Value = 27657.642
$ws.Range("N132").Value = -32717.642

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1069.0714
$ws.Range("I64").Value = 512.25
$ws.Range("K64").Value = 512.25
$ws.Range("M64").Value = -287.25
$ws.Range("H67").Value = 1069.0714
$ws.Range("I67").Value = 512.25
$ws.Range("K67").Value = 512.25
$ws.Range("M67").Value = 267.75
$ws.Range("H94").Value = 1372.3793
$ws.Range("J94").Value = 1174.0834
$ws.Range("L94").Value = 1174.0834
$ws.Range("N94").Value = -2076.0834
$ws.Range("H134").Value = 7165.25
$ws.Range("I134").Value = 5704.9287
$ws.Range("K134").Value = 17114.7861
$ws.Range("M134").Value = -14579.7861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5749.5
$ws.Range("I31").Value = 3138.5557
$ws.Range("K31").Value = 3138.5557
$ws.Range("M31").Value = -2843.5557
$ws.Range("H34").Value = 5749.5
$ws.Range("I34").Value = 3138.5557
$ws.Range("K34").Value = 3138.5557
$ws.Range("M34").Value = -2936.5557
$ws.Range("H132").Value = 17930.13
$ws.Range("I132").Value = 1813.2858
$ws.Range("J132").Value = 63057.3
$ws.Range("K132").Value = 5439.857400000001
$ws.Range("L132").Value = 189171.9
$ws.Range("M132").Value = -2909.857400000001
$ws.Range("N132").Value = -194231.9
$ws.Range("H134").Value = 5973
$ws.Range("I134").Value = 6178.4
$ws.Range("J134").Value = 5630.6665
$ws.Range("K134").Value = 18535.2
$ws.Range("L134").Value = 16891.9995
$ws.Range("M134").Value = -16000.2
$ws.Range("N134").Value = -21961.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1149
$ws.Range("J97").Value = 1250
$ws.Range("L97").Value = 1250
$ws.Range("N97").Value = -2242
$ws.Range("H113").Value = 1569.75
$ws.Range("I113").Value = 1569.75
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1569.75
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 600.25
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 6924.0435
$ws.Range("I132").Value = 4646.6665
$ws.Range("K132").Value = 13939.9995
$ws.Range("M132").Value = -11409.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 134.83333
$ws.Range("I55").Value = 113.46154
$ws.Range("K55").Value = 113.46154
$ws.Range("M55").Value = 59.53846
$ws.Range("H56").Value = 12409
$ws.Range("I56").Value = 4914.4
$ws.Range("J56").Value = 24900
$ws.Range("K56").Value = 4914.4
$ws.Range("L56").Value = 24900
$ws.Range("M56").Value = -4223.4
$ws.Range("N56").Value = -26282
$ws.Range("H82").Value = 1587
$ws.Range("I82").Value = 1510.4445
$ws.Range("J82").Value = 1759.25
$ws.Range("K82").Value = 1510.4445
$ws.Range("L82").Value = 1759.25
$ws.Range("M82").Value = -1149.4445
$ws.Range("N82").Value = -2481.25
$ws.Range("H85").Value = 1587
$ws.Range("I85").Value = 1510.4445
$ws.Range("J85").Value = 1759.25
$ws.Range("K85").Value = 1510.4445
$ws.Range("L85").Value = 1759.25
$ws.Range("M85").Value = -262.4445000000001
$ws.Range("N85").Value = -4255.25
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H93").Value = 1225.25
$ws.Range("I93").Value = 1134
$ws.Range("J93").Value = 1499
$ws.Range("K93").Value = 1134
$ws.Range("L93").Value = 1499
$ws.Range("M93").Value = 114
$ws.Range("N93").Value = -3995
$ws.Range("H132").Value = 3859.7908
$ws.Range("I132").Value = 2470.1667
$ws.Range("K132").Value = 7410.500100000001
$ws.Range("M132").Value = -4880.500100000001
$ws.Range("H134").Value = 82714.5
$ws.Range("J134").Value = 82714.5
$ws.Range("L134").Value = 82714.5
$ws.Range("N134").Value = -92854.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5058.3335
$ws.Range("J81").Value = 7474.5
$ws.Range("L81").Value = 14949
$ws.Range("N81").Value = -17071
$ws.Range("H84").Value = 5058.3335
$ws.Range("J84").Value = 7474.5
$ws.Range("L84").Value = 74745
$ws.Range("N84").Value = -85353
$ws.Range("H100").Value = 1771.4117
$ws.Range("I100").Value = 1781.3
$ws.Range("J100").Value = 1757.2858
$ws.Range("K100").Value = 3562.6
$ws.Range("L100").Value = 3514.5716
$ws.Range("M100").Value = -3021.6
$ws.Range("N100").Value = -4596.5716
$ws.Range("H107").Value = 1091.0731
$ws.Range("I107").Value = 769.2121
$ws.Range("J107").Value = 2418.75
$ws.Range("K107").Value = 2307.6363
$ws.Range("L107").Value = 7256.25
$ws.Range("M107").Value = -387.6363000000001
$ws.Range("N107").Value = -11096.25
